$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.476.33"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "2.668.78"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.33"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.37"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "2.668.82"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.20"
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.80"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D17").Value = "67.413.35"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "2.669.72"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.70"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.83"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.76"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.78"
$ws.Range("E23").Value = "  -4.01%  "
$ws.Range("E24").Value = "  -5.72%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.73"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "2.813.78"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  -4.53%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.52"
$ws.Range("E31").Value = "  -8.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.43"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.50"
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.92"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -8.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.31"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "0.0₆0299"
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.51"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  -5.03%  "
